$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on D (price) and E (volume%) columns so Excel
# does not reinterpret numeric-looking strings (e.g. "5.38", "1.00") as numbers,
# matching the source data which stores these as plain text (inline strings).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "62.859.51"
$ws.Range("E2").Value = "  -1.86%  "
$ws.Range("D3").Value = "3.229.82"
$ws.Range("E3").Value = "  -2.56%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").Value = "594.84"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").Value = "136.28"
$ws.Range("E6").Value = "  -5.48%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").Value = "3.230.60"
$ws.Range("E8").Value = "  -2.43%  "
$ws.Range("E9").Value = "  -2.67%  "
$ws.Range("E10").Value = "  -3.11%  "
$ws.Range("D11").Value = "5.38"
$ws.Range("E11").Value = "  -1.63%  "
$ws.Range("E12").Value = "  -3.70%  "
$ws.Range("D13").Value = "0.0000240"
$ws.Range("E13").Value = "  -3.94%  "
$ws.Range("D14").Value = "33.69"
$ws.Range("E14").Value = "  -4.03%  "
$ws.Range("D15").Value = "3.763.86"
$ws.Range("E15").Value = "  -1.82%  "
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").Value = "3.231.83"
$ws.Range("E17").Value = "  -1.76%  "
$ws.Range("D18").Value = "62.945.49"
$ws.Range("E18").Value = "  -1.50%  "
$ws.Range("D19").Value = "6.72"
$ws.Range("E19").Value = "  -3.17%  "
$ws.Range("D20").Value = "468.59"
$ws.Range("E20").Value = "  -3.36%  "
$ws.Range("E21").Value = "  -3.90%  "
$ws.Range("D22").Value = "0.718"
$ws.Range("E22").Value = "  -3.85%  "
$ws.Range("D23").Value = "7.73"
$ws.Range("E23").Value = "  -4.50%  "
$ws.Range("D24").Value = "13.49"
$ws.Range("E24").Value = "  -1.10%  "
$ws.Range("D25").Value = "84.49"
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").Value = "2.70"
$ws.Range("E27").Value = "  -2.51%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("D29").Value = "7.92"
$ws.Range("E29").Value = "  -4.57%  "
$ws.Range("D30").Value = "6.92"
$ws.Range("E30").Value = "  -5.36%  "
$ws.Range("D31").Value = "2.08"
$ws.Range("E31").Value = "  -4.50%  "
$ws.Range("D32").Value = "27.92"
$ws.Range("E32").Value = "  -1.43%  "
$ws.Range("E33").Value = "  -6.08%  "
$ws.Range("E34").Value = "  -5.34%  "
$ws.Range("E35").Value = "  -4.20%  "
$ws.Range("D36").Value = "5.88"
$ws.Range("E36").Value = "  -2.80%  "
$ws.Range("D37").Value = "51.78"
$ws.Range("E37").Value = "  -2.58%  "
$ws.Range("D38").Value = "0.0₃0712"
$ws.Range("E38").Value = "  -4.38%  "
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "421.66"
$ws.Range("E40").Value = "  -1.83%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "3.035.25"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D42").Value = "0.117"
$ws.Range("E42").Value = "  +5.44%  "
$ws.Range("D43").Value = "8.12"
$ws.Range("E43").Value = "  -4.53%  "
$ws.Range("D44").Value = "2.63"
$ws.Range("E44").Value = "  -6.44%  "
$ws.Range("E45").Value = "  -5.87%  "
$ws.Range("E46").Value = "  -4.68%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "126.47"
$ws.Range("E48").Value = "  +2.45%  "
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").Value = "35.57"
$ws.Range("E49").Value = "  +3.71%  "
$ws.Range("D50").Value = "25.70"
$ws.Range("E50").Value = "  -2.78%  "
$ws.Range("E51").Value = "  -2.57%  "

# Restore default (Normal) style on the touched price/volume cells so the
# cell styling matches the original (no explicit style index).
$ws.Range("D2:E51").Style = "Normal"
